$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125, shifting existing rows 125:143 down to 126:144
$ws.Rows("125:125").Insert()

# Populate the newly inserted row 125 with the new data record
$ws.Cells.Item(125,1).Value  = 11
$ws.Cells.Item(125,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(125,3).Value  = "Bíobío"
$ws.Cells.Item(125,4).Value  = 44505
$ws.Cells.Item(125,4).NumberFormat = $ws.Cells.Item(126,4).NumberFormat
$ws.Cells.Item(125,5).Value  = 8
$ws.Cells.Item(125,6).Value  = 100114013
$ws.Cells.Item(125,7).Value  = "Zanahoria"
$ws.Cells.Item(125,8).Value  = "Sin especificar"
$ws.Cells.Item(125,9).Value  = "Primera"
$ws.Cells.Item(125,10).Value = 250
$ws.Cells.Item(125,11).Value = 8000
$ws.Cells.Item(125,12).Value = 8500
$ws.Cells.Item(125,13).Value = 8200
$ws.Cells.Item(125,14).Value = "$/saco 20 kilos"
$ws.Cells.Item(125,15).Value = "Provincia del Elquí"
$ws.Cells.Item(125,16).Value = 410
$ws.Cells.Item(125,17).Value = 20
$ws.Cells.Item(125,18).Value = "Hortaliza"
